# Weekly update: a new price record for "Apio" (Vega Monumental Concepción)
# is inserted as a new row 133, pushing the existing rows 133:149 down to
# 134:150 (dimension grows from A1:R149 to A1:R150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 133; this shifts rows 133-149 -> 134-150
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly record.
$ws.Cells.Item(133, 1).Value  = 11
$ws.Cells.Item(133, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(133, 3).Value  = "Bíobío"
$ws.Cells.Item(133, 4).Value  = 44476
$ws.Cells.Item(133, 5).Value  = 8
$ws.Cells.Item(133, 6).Value  = 100112017
$ws.Cells.Item(133, 7).Value  = "Apio"
$ws.Cells.Item(133, 8).Value  = "Americana (o)"
$ws.Cells.Item(133, 9).Value  = "Primera"
$ws.Cells.Item(133, 10).Value = 150
$ws.Cells.Item(133, 11).Value = 7500
$ws.Cells.Item(133, 12).Value = 8500
$ws.Cells.Item(133, 13).Value = 8000
$ws.Cells.Item(133, 14).Value = "`$/docena de matas"
$ws.Cells.Item(133, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(133, 16).Value = 1333
$ws.Cells.Item(133, 17).Value = 6
$ws.Cells.Item(133, 18).Value = "Hortaliza"
